$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename Resistance header, add new Torque header, move the
#     denominator constant from I1 to M1 (frees up the old I1 cell). ---
$ws.Range("H1").Value = "Tm"
$ws.Range("I1").ClearContents()
$ws.Range("L1").Value = "Rr"
$ws.Range("M1").Value = 5.3

# --- Row 2: nothing changes in A/B/C, but we now add a Torque formula.
#     B2 is still the text "-" and F2 is blank, so B2*F2 correctly errors
#     out with #VALUE!, matching the source data. ---
$ws.Range("H2").Formula = "=B2*F2"

# --- Rows 3-6: the "Ia" (now "R") column was previously stored as text
#     labels pulled from the shared-string table; the new data makes them
#     real numeric measurements. ---
$ws.Range("B3").Value = 2.18
$ws.Range("B4").Value = 3.08
$ws.Range("B5").Value = 4.32
$ws.Range("B6").Value = 5.4

# --- F column formulas now divide by $M$1 instead of $I$1 (since the
#     constant moved). Re-enter each one individually. ---
$ws.Range("F3").Formula = "=E3/`$M`$1"
$ws.Range("F4").Formula = "=E4/`$M`$1"
$ws.Range("F5").Formula = "=E5/`$M`$1"
$ws.Range("F6").Formula = "=E6/`$M`$1"

# --- New Torque column: H = B * F, entered once across the block so Excel
#     keeps it as a single shared formula. ---
$ws.Range("H3:H6").Formula = "=B3*F3"
$ws.Range("H3:H6").NumberFormat = "0.000"

# --- Selection moves to the last-entered cell. ---
$ws.Range("H6").Select()
